$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent results for the 380 kV case (rows 2-25)
$data = @(
    @(2, 14.63600841952538, 4.562363476012342, 8.577642651360252, 10.07149674727916, 43.39203069139801, 35.38290916615532, 12.91868258295429, 10.38218828623057),
    @(3, 14.53035442462144, 4.327100145637402, 8.569480945865768, 10.07816099576439, 43.04555017766513, 35.21440922354945, 12.84548908304106, 10.37515367810482),
    @(4, 14.47017483856933, 4.174976701742761, 8.564290945351265, 10.08373185145526, 42.83813177391508, 35.114732043346, 12.80449573145657, 10.37281800431302),
    @(5, 14.44685833451983, 4.111083119008001, 8.562130747607892, 10.08637415429245, 42.75500481761795, 35.0750823571013, 12.78879970875627, 10.37236603812561),
    @(6, 14.44306030516648, 4.100359658480728, 8.561769307758373, 10.0868353877573, 42.74128775210512, 35.06855763001295, 12.78625478861138, 10.37232120366007),
    @(7, 14.4698554621125, 4.174122672118275, 8.564261995716487, 10.08376597947053, 42.83700495902774, 35.11419336476499, 12.80427994323951, 10.37280988393503),
    @(8, 14.59862207325103, 4.482851133118112, 8.574865400116463, 10.07348778004934, 43.27149609731149, 35.32403676911838, 12.89263647667355, 10.37935176427853),
    @(9, 14.88705600075642, 5.026492430601511, 8.594254565786688, 10.06505484436644, 44.16303208732251, 35.7646908613316, 13.09642447285619, 10.40786274426249),
    @(10, 15.11897000663661, 5.387340573114787, 8.60766294976967, 10.06598543810532, 44.83832896536194, 36.10499053265389, 13.2635700373845, 10.4382786045772),
    @(11, 15.22837739660818, 5.543005406550313, 8.613585109635419, 10.06795009638274, 45.14910073752393, 36.26311309547302, 13.34311096286081, 10.45414622569842),
    @(12, 15.27033045246955, 5.600724470006026, 8.615802476044678, 10.06891499163132, 45.26722323204427, 36.32344156677219, 13.37370981518668, 10.46044431305018),
    @(13, 15.26127249077526, 5.588348340289751, 8.615326044066638, 10.06869737113998, 45.24176506018661, 36.31042912574032, 13.36709894507525, 10.45907508433289),
    @(14, 15.23181864872218, 5.547778630995938, 8.61376803769976, 10.06802505576878, 45.15881035606585, 36.26806747242104, 13.34561891073661, 10.45465859054441),
    @(15, 15.21384423198961, 5.522768429791227, 8.612810436764072, 10.06764199067733, 45.10805328720986, 36.24217765986629, 13.33252329740144, 10.45199095818195),
    @(16, 15.11189529780674, 5.376996046592692, 8.60727237736749, 10.06588799932032, 44.81808492476359, 36.09472139953011, 13.25844021986137, 10.4372822659629),
    @(17, 15.05032642878134, 5.285389715672614, 8.603829726255775, 10.0652062733346, 44.64106200681015, 36.00509327827845, 13.21387256300066, 10.42877736986076),
    @(18, 15.01528395070779, 5.231901645928957, 8.601832891910272, 10.06495928893035, 44.53958418092962, 35.95385554367749, 13.18856998981828, 10.42407695578419),
    @(19, 15.00348394855172, 5.213654620134081, 8.601153919124748, 10.06490061010264, 44.50528636231941, 35.93656204238652, 13.18006068579738, 10.42251842409068),
    @(20, 15.05684249563254, 5.295224071800511, 8.604197929994092, 10.06526382910189, 44.65987163447039, 36.01460203600135, 13.21858273346689, 10.42966294265955),
    @(21, 15.24045608285896, 5.559728309406314, 8.614226344374293, 10.06821654156548, 45.18316481762319, 36.2804980938411, 13.35191534222235, 10.45594799180333),
    @(22, 15.36348700299612, 5.725438409135957, 8.620633578023527, 10.07143360662583, 45.5276968361673, 36.45689187034017, 13.44183203980066, 10.47481195257577),
    @(23, 15.29755911430258, 5.637652738417959, 8.617227267956816, 10.06959907422474, 45.34360681924166, 36.3625167995848, 13.39359650206069, 10.46459068685408),
    @(24, 15.05389547528822, 5.29078052035626, 8.604031519944485, 10.0652373565995, 44.65136687718072, 36.01030221919768, 13.21645226680449, 10.4292619857801),
    @(25, 14.80536995359458, 4.886144561423396, 8.589157607797157, 10.06608289690151, 43.91800344213554, 35.6424827661212, 13.03814341453408, 10.39847915149011)
)

$colIndexes = @(2, 3, 4, 5, 6, 9, 11, 12)  # B, C, D, E, F, I, K, L

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $c = $colIndexes[$i]
        $value = $entry[$i + 1]
        $ws.Cells.Item($r, $c).Value = $value
    }
}

Write-Host "Updated $($data.Length) rows x $($colIndexes.Length) columns of loading_percent data"
